$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the no-op left-alignment style from header/unit rows (C1:Q1, C3:Q3)
$ws.Range("C1:Q1").Style = "Normal"
$ws.Range("C3:Q3").Style = "Normal"

# Add new "grade" column
$ws.Range("R1").Value = "grade"

$ws.Range("R4").Value = "D"
$ws.Range("R5").Value = "C"
$ws.Range("R6").Value = "F"
$ws.Range("R7").Value = "F"
$ws.Range("R8").Value = "F"
$ws.Range("R9").Value = "F"
$ws.Range("R10").Value = "F"
$ws.Range("R11").Value = "F"
$ws.Range("R12").Value = "E"
$ws.Range("R13").Value = "F"
$ws.Range("R14").Value = "F"
$ws.Range("R15").Value = "A+"
$ws.Range("R16").Value = "C"
$ws.Range("R17").Value = "A+"
$ws.Range("R18").Value = "A+"
$ws.Range("R19").Value = "D"
$ws.Range("R20").Value = "D"
$ws.Range("R21").Value = "D"
$ws.Range("R22").Value = "E"
$ws.Range("R23").Value = "F"
$ws.Range("R24").Value = "F"
$ws.Range("R25").Value = "F"
$ws.Range("R26").Value = "E"
$ws.Range("R27").Value = "D"
$ws.Range("R28").Value = "E"
$ws.Range("R29").Value = "C"
$ws.Range("R30").Value = "E"
$ws.Range("R31").Value = "B"
$ws.Range("R32").Value = "D"
$ws.Range("R33").Value = "D"
$ws.Range("R34").Value = "C"
$ws.Range("R35").Value = "E"
$ws.Range("R36").Value = "E"
$ws.Range("R37").Value = "B"
$ws.Range("R38").Value = "D"
$ws.Range("R39").Value = "F"
$ws.Range("R40").Value = "F"
$ws.Range("R41").Value = "E"
$ws.Range("R42").Value = "D"
$ws.Range("R43").Value = "B"
$ws.Range("R44").Value = "C"
$ws.Range("R45").Value = "F"
$ws.Range("R46").Value = "E"
$ws.Range("R47").Value = "C"
$ws.Range("R48").Value = "F"
$ws.Range("R49").Value = "C"
$ws.Range("R50").Value = "B"
$ws.Range("R51").Value = "B"
$ws.Range("R52").Value = "A"
$ws.Range("R53").Value = "A+"
$ws.Range("R54").Value = "A+"
$ws.Range("R55").Value = "D"
$ws.Range("R56").Value = "C"
$ws.Range("R57").Value = "B"
$ws.Range("R58").Value = "B"
$ws.Range("R59").Value = "F"
$ws.Range("R60").Value = "D"
$ws.Range("R61").Value = "F"
$ws.Range("R62").Value = "D"
$ws.Range("R63").Value = "B"
$ws.Range("R64").Value = "B"
$ws.Range("R65").Value = "E"
$ws.Range("R66").Value = "B"
$ws.Range("R67").Value = "D"
$ws.Range("R68").Value = "B"
$ws.Range("R69").Value = "C"
$ws.Range("R70").Value = "B"
$ws.Range("R71").Value = "E"
$ws.Range("R72").Value = "A+"

# Match the workbook's final view/selection state
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R2:R3").Select() | Out-Null
